$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RecipeProbability")

$ws.Range("C6").Value = 19.0
$ws.Range("C7").Value = 1.0
$ws.Range("C9").Value = 19.0
$ws.Range("C11").Value = 1.0
$ws.Range("C12").Value = 19.0
$ws.Range("C13").Value = 5.0
$ws.Range("C14").Value = 7.0
